# Update countries & provincias Spain
# - reorder "Suazilandia"/"Mayotte" (rows 117/118) and
#   "Santa Lucia"/"Timor Oriental" (rows 202/203)
# - refresh the "Datos actualizados..." timestamp string
# - refresh the day's numeric counters for a handful of countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country name order (A column) --------------------------------
$ws.Range("A117").Value = "Suazilandia"
$ws.Range("A118").Value = "Mayotte"

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Timestamp string ----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 20:47"

# --- Numeric refresh: Estados Unidos (row 4) -----------------------------
$ws.Range("B4").Value = 5120955
$ws.Range("C4").Value = 25431
$ws.Range("D4").Value = 2620147
$ws.Range("E4").Value = 2336231
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 483
$ws.Range("H4").Value = 164577

# --- Numeric refresh: India (row 6) --------------------------------------
$ws.Range("B6").Value = 2152020
$ws.Range("C6").Value = 65156
$ws.Range("D6").Value = 1479804
$ws.Range("E6").Value = 628763
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 875
$ws.Range("H6").Value = 43453

# --- Numeric refresh: Turquia (row 20) -----------------------------------
$ws.Range("B20").Value = 239622
$ws.Range("C20").Value = 1172
$ws.Range("D20").Value = 222656
$ws.Range("E20").Value = 11137
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 5829

# --- Numeric refresh: Alemania (row 22) ----------------------------------
$ws.Range("B22").Value = 216692
$ws.Range("C22").Value = 377
$ws.Range("D22").Value = 197400
$ws.Range("E22").Value = 10033
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 9259

# --- Numeric refresh: Republica Dominicana (row 38) ----------------------
$ws.Range("B38").Value = 78778
$ws.Range("C38").Value = 1069
$ws.Range("D38").Value = 42538
$ws.Range("E38").Value = 34951
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 30
$ws.Range("H38").Value = 1289

# --- Numeric refresh: Libano (row 100) -----------------------------------
$ws.Range("B100").Value = 6223
$ws.Range("C100").Value = 272
$ws.Range("D100").Value = 2043
$ws.Range("E100").Value = 4102
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 8
$ws.Range("H100").Value = 78

# --- Numeric refresh: Malaui (row 109) -----------------------------------
$ws.Range("B109").Value = 4624
$ws.Range("C109").Value = 49
$ws.Range("D109").Value = 2329
$ws.Range("E109").Value = 2152
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 6
$ws.Range("H109").Value = 143

# --- Numeric refresh: row 117 (now Suazilandia) --------------------------
$ws.Range("B117").Value = 3128
$ws.Range("C117").Value = 92
$ws.Range("D117").Value = 1565
$ws.Range("E117").Value = 1507
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 56

# --- Numeric refresh: row 118 (now Mayotte) ------------------------------
$ws.Range("B118").Value = 3068
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 2835
$ws.Range("E118").Value = 194
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 39

# --- Numeric refresh: Burkina Faso (row 147) -----------------------------
$ws.Range("B147").Value = 1175
$ws.Range("C147").Value = 17
$ws.Range("D147").Value = 974
$ws.Range("E147").Value = 147
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 54
